$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.234.33"
$ws.Range("E2").Value = "  -2.59%  "
$ws.Range("D3").Value = "2.569.37"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("D13").Value = "3.024.61"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").Value = "58.141.99"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").Value = "2.565.84"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.418"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "0.0₃0727"
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.848"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("E36").Value = "  -5.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.811"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("E38").Value = "  -3.98%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "277.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.586"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.11%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "1.900.62"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
